$d = $word.ActiveDocument

# ======================================================================
# Helper: build a zero-length Range at an absolute character offset
# inside a paragraph, without suffering from the stale-cache issue that
# $d.Range(pos, pos) exhibits right after the document has grown.
# ======================================================================
function PointRangeFromParaStart($para, $paraStart, $offset) {
    $rr = $para.Range.Duplicate
    $rr.Collapse(1)
    $rr.MoveStart(1, $offset - $paraStart)
    $rr.Collapse(1)
    return $rr
}

# ----------------------------------------------------------------------
# 1. Fix the split "Alison J. F" + bookmark + "rancis" into a single run
#    " Alison J. Francis" (the _GoBack bookmark that used to sit between
#    "F" and "rancis" is removed from here; it gets relocated below,
#    next to the newly added "Reference: ..." paragraph).
# ----------------------------------------------------------------------

$find1 = $d.Content.Find
$find1.Execute(" Alison J. Francis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetRange = $find1.Parent.Duplicate
$startPos = $targetRange.Start
$endPos = $targetRange.End

# Temporary bookmarks act as merge barriers so that only the text inside
# the target span coalesces into a single run; the surrounding runs
# ("," before, and "(2015)..." after) are left completely untouched.
$d.Bookmarks.Add("ZZTempBarrierBefore", $d.Range($startPos, $startPos))
$d.Bookmarks.Add("ZZTempBarrierAfter", $d.Range($endPos, $endPos))

# Remove the original _GoBack bookmark that sits between "F" and "rancis".
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# A genuine text change (to a placeholder) forces the run coalescing...
$find2 = $d.Content.Find
$find2.Execute(" Alison J. Francis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find2.Parent.Text = " Alison J. FrancisZZPLACEHOLDERZZ"

# ...then a second genuine edit removes the placeholder, leaving the
# merged run with the desired final text, without re-splitting it.
$find3 = $d.Content.Find
$find3.Execute(" Alison J. FrancisZZPLACEHOLDERZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find3.Parent.Text = " Alison J. Francis"

$d.Bookmarks.Item("ZZTempBarrierBefore").Delete()
$d.Bookmarks.Item("ZZTempBarrierAfter").Delete()

# ----------------------------------------------------------------------
# 2. Append two new paragraphs at the end of the document: a
#    "Reference: ..." paragraph (re-using the same citation, with the
#    _GoBack bookmark now anchored right before "(2015)") and a
#    "Project Management Journal..." citation paragraph.
# ----------------------------------------------------------------------

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$lastPara.Range.InsertParagraphAfter()

$refPara = $d.Paragraphs.Item($lastParaIndex + 1)
$refRange = $refPara.Range
$refRange.Collapse(1)
$paraStart = $refRange.Start

$seg1 = "Reference: "
$seg2 = "Alicia Medina"
$seg3 = ","
$seg4 = " Alison J. Francis"
$seg5 = " "
$seg6 = "(2015)  "
$fullText = $seg1 + $seg2 + $seg3 + $seg4 + $seg5 + $seg6

$refRange.InsertAfter($fullText)

$p1 = $paraStart + $seg1.Length
$p2 = $p1 + $seg2.Length
$p3 = $p2 + $seg3.Length
$p4 = $p3 + $seg4.Length
$p5 = $p4 + $seg5.Length

# Split the single freshly-inserted run into the six separate runs shown
# in the target diff by dropping temporary bookmarks at each boundary
# (inserting a bookmark splits the run it lands inside of), then remove
# the temporary ones, leaving only _GoBack in its final resting place.
$d.Bookmarks.Add("ZZBarrier1", (PointRangeFromParaStart $refPara $paraStart $p1))
$d.Bookmarks.Add("ZZBarrier2", (PointRangeFromParaStart $refPara $paraStart $p2))
$d.Bookmarks.Add("ZZBarrier3", (PointRangeFromParaStart $refPara $paraStart $p3))
$d.Bookmarks.Add("ZZBarrier4", (PointRangeFromParaStart $refPara $paraStart $p4))
$d.Bookmarks.Add("_GoBack", (PointRangeFromParaStart $refPara $paraStart $p5))

$d.Bookmarks.Item("ZZBarrier1").Delete()
$d.Bookmarks.Item("ZZBarrier2").Delete()
$d.Bookmarks.Item("ZZBarrier3").Delete()
$d.Bookmarks.Item("ZZBarrier4").Delete()

$refPara.Range.InsertParagraphAfter()

$citePara = $d.Paragraphs.Item($lastParaIndex + 2)
$citeRange = $citePara.Range
$citeRange.Collapse(1)
$citeRange.InsertAfter("Project Management Journal, Vol. 46, No. 5, 81" + [char]0x2013 + "93 " + [char]0x00A9 + " 2015 by the Project Management Institute Published online in Wiley Online Library (wileyonlinelibrary.com). DOI: 10.1002/pmj.21530")

Write-Host "Done."
